$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are stored as literal text in the source data (not numbers),
# e.g. '27.797.33' or '1.000' with significant trailing zeros. Force text format on the
# D column for the rows we touch so Excel does not auto-convert these into numbers,
# which would silently corrupt values like '1.000' -> 1 or '0.00001050' -> 1.05E-05.
# NumberFormat is switched back to General immediately afterwards so the cell's
# displayed/stored formatting is unaffected - only the literal text content changes.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.763.22"
$ws.Range("E2").Value = "  -1.43%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.750.44"
$ws.Range("E3").Value = "  -3.72%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.82"
$ws.Range("E5").Value = "  -2.60%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.31%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4218"
$ws.Range("E7").Value = "  -4.54%  "

# Row 8
$ws.Range("E8").Value = "  -3.08%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.39"
$ws.Range("E9").Value = "  -5.41%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07431"
$ws.Range("E10").Value = "  -3.68%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.083"
$ws.Range("E11").Value = "  -2.97%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.45%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.68"
$ws.Range("E13").Value = "  -5.96%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.033"
$ws.Range("E14").Value = "  -4.34%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.262"
$ws.Range("E15").Value = "  -3.47%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.763.76"
$ws.Range("E16").Value = "  -3.59%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.73"
$ws.Range("E17").Value = "  -2.57%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001050"
$ws.Range("E18").Value = "  -2.94%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06343"
$ws.Range("E19").Value = "  -2.10%  "

# Row 20
$ws.Range("E20").Value = "  -0.24%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.99"
$ws.Range("E21").Value = "  -3.23%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.903"
$ws.Range("E22").Value = "  -6.25%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.783.08"
$ws.Range("E23").Value = "  -1.62%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.13"
$ws.Range("E24").Value = "  -4.79%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.100"
$ws.Range("E25").Value = "  +0.06%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.02"
$ws.Range("E26").Value = "  +1.12%  "

# Row 27
$ws.Range("E27").Value = "  -2.69%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.963.75"
$ws.Range("E28").Value = "  -3.41%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.118"
$ws.Range("E29").Value = "  -9.33%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.51"
$ws.Range("E30").Value = "  -3.49%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.116"
$ws.Range("E31").Value = "  -6.48%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.640"

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.539"
$ws.Range("E33").Value = "  -5.29%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08838"
$ws.Range("E34").Value = "  -4.60%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.23"
$ws.Range("E35").Value = "  -5.67%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02279"
$ws.Range("E36").Value = "  -2.12%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6277"
$ws.Range("E39").Value = "  -4.25%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.920"
$ws.Range("E40").Value = "  -4.52%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.172"
$ws.Range("E41").Value = "  -2.30%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9994"
$ws.Range("E42").Value = "  -0.38%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.394"
$ws.Range("E43").Value = "  +0.37%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.832"
$ws.Range("E44").Value = "  -3.45%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.33"
$ws.Range("E45").Value = "  -4.05%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5845"
$ws.Range("E46").Value = "  -3.58%  "

# Row 47
$ws.Range("E47").Value = "  -2.44%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.30"
$ws.Range("E48").Value = "  -3.40%  "

# Row 49
$ws.Range("E49").Value = "  -3.67%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.176"
$ws.Range("E50").Value = "  +2.24%  "

# Row 51
$ws.Range("E51").Value = "  -2.76%  "

# Rows 37 and 38 swap content (Algorand <-> Hedera) and both get updated D/E values.
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06018"
$ws.Range("E37").Value = "  -2.31%  "

$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2088"
$ws.Range("E38").Value = "  -4.03%  "
